# Apply cell-value updates described by the commit diff (prices / 1h volume %,
# a block of Coin/Link rows that got re-sorted, and a couple of refreshed % figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) store numeric-looking text (e.g. "255.69", "4.19%").
# A leading apostrophe forces Excel to keep them as literal text instead of coercing to
# a number, matching how the workbook already stores these columns.

# Row 2
$ws.Range("D2").Value = "`'255.69"
$ws.Range("E2").Value = "`'4.19%"
# Row 3
$ws.Range("D3").Value = "`'27.93"
$ws.Range("E3").Value = "`'-4.82%"
# Row 4
$ws.Range("D4").Value = "`'5.350"
$ws.Range("E4").Value = "`'3.89%"
# Row 5
$ws.Range("D5").Value = "`'0.05829"
$ws.Range("E5").Value = "`'0.94%"
# Row 6
$ws.Range("D6").Value = "`'6.711"
$ws.Range("E6").Value = "`'1.23%"
# Row 7
$ws.Range("E7").Value = "`'2.24%"
# Row 8
$ws.Range("D8").Value = "`'0.8711"
$ws.Range("E8").Value = "`'1.72%"
# Row 9
$ws.Range("D9").Value = "`'0.9116"
$ws.Range("E9").Value = "`'6.66%"
# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "`'0.1421"
$ws.Range("E10").Value = "`'4.06%"
# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "`'0.07224"
$ws.Range("E11").Value = "`'2.22%"
# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "`'0.03180"
$ws.Range("E12").Value = "`'3.93%"
# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "`'0.09242"
$ws.Range("E13").Value = "`'-1.39%"
# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "`'0.001542"
$ws.Range("E14").Value = "`'1.33%"
# Row 15
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "`'0.0006046"
$ws.Range("E15").Value = "`'0.77%"
# Row 16
$ws.Range("D16").Value = "`'0.005965"
$ws.Range("E16").Value = "`'-0.55%"
# Row 17
$ws.Range("E17").Value = "`'0.40%"
# Row 18
$ws.Range("E18").Value = "`'5.01%"
# Row 19
$ws.Range("D19").Value = "`'0.3167"
$ws.Range("E19").Value = "`'-1.13%"
# Row 21
$ws.Range("E21").Value = "`'2.04%"
# Row 22
$ws.Range("D22").Value = "`'3.524"
$ws.Range("E22").Value = "`'10.83%"
# Row 23
$ws.Range("D23").Value = "`'0.04153"
$ws.Range("E23").Value = "`'0.57%"
# Row 25
$ws.Range("D25").Value = "`'0.001223"
$ws.Range("E25").Value = "`'-0.28%"
# Row 26
$ws.Range("D26").Value = "`'0.004874"
$ws.Range("E26").Value = "`'18.00%"
# Row 27
$ws.Range("E27").Value = "`'-0.90%"
# Row 28
$ws.Range("E28").Value = "`'0.62%"
# Row 40
$ws.Range("D40").Value = "`'0.03852"
$ws.Range("E40").Value = "`'3.38%"
# Row 41
$ws.Range("D41").Value = "`'0.005768"
$ws.Range("E41").Value = "`'64.84%"
# Row 42
$ws.Range("D42").Value = "`'0.1102"
$ws.Range("E42").Value = "`'2.98%"
# Row 43
$ws.Range("E43").Value = "`'-9.90%"
# Row 44
$ws.Range("D44").Value = "`'0.009927"
$ws.Range("E44").Value = "`'8.47%"
# Row 45
$ws.Range("D45").Value = "`'0.00005268"
$ws.Range("E45").Value = "`'-0.26%"
# Row 46
$ws.Range("E46").Value = "`'-0.07%"
# Row 47
$ws.Range("D47").Value = "`'0.09993"
# Row 48
$ws.Range("D48").Value = "`'0.002135"
$ws.Range("E48").Value = "`'-1.76%"
# Row 49
$ws.Range("E49").Value = "`'-0.07%"
# Row 50
$ws.Range("D50").Value = "`'0.0001999"
$ws.Range("E50").Value = "`'-0.07%"
